$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Append four new Box product links below the existing one (A2), continuing
# the list started in A1 (header) / A2 (first link).
$ws.Range("A3").Value = "https://box.co.uk/nx-j2kek-004-acer-swift-14-intel-core-ultra-7-256v-1"
$ws.Range("A4").Value = "https://box.co.uk/nx-kyxek-003-acer-swift-14-qualcomm-snapdragon-100"
$ws.Range("A5").Value = "https://box.co.uk/nx-j2kek-002-acer-swift-14-intel-core-ultra-5-226v-1"
$ws.Range("A6").Value = "https://box.co.uk/nx-ab1ek-00d-acer-swift-3-amd-ryzen-5-5500u-16-gb"

# Match the visual "Hyperlink" look already used for A2 on the new rows, and
# also on the trailing blank row (A7) left over from the paste.
$ws.Range("A3:A7").Style = "Hyperlink"

# Leave the selection where the author's cursor ended up after pasting.
[void]$ws.Range("A8").Select()
